$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: remove the existing "_GoBack" bookmark; we'll re-create it
# later at the end of the "Readme" paragraph.
# ------------------------------------------------------------------
$gb = $d.Bookmarks.Item("_GoBack")
$gb.Delete()

# ------------------------------------------------------------------
# Step 2: locate the "By: hxr190001, pxr180025." paragraph and turn it
# into "By: Harshita Rastogi (hxr190001)," followed by a new paragraph
# "Prajakta Ray (pxr180025)."
# ------------------------------------------------------------------
$byPara = $d.Paragraphs.Item(2)
$byRange = $byPara.Range

# Insert ")" right after "hxr190001" (before the ", pxr180025." tail).
$hxrRange = $d.Range($byRange.Start, $byRange.End)
$hxrRange.Find.Execute("hxr190001") | Out-Null
$hxrRange.Collapse(0)             # collapse to end of "hxr190001"
$hxrRange.InsertBefore(")")

# Insert "Harshita Rastogi (" right before "hxr190001".
$byRange2 = $d.Range($byPara.Range.Start, $byPara.Range.End)
$byRange2.Find.Execute("hxr190001") | Out-Null
$byRange2.Collapse(1)              # collapse to start of "hxr190001"
$byRange2.InsertBefore("Harshita Rastogi (")

# Split the paragraph right after the comma that follows ")" -- i.e.
# turn "...hxr190001), pxr180025." into "...hxr190001)," + new
# paragraph " pxr180025."
$byRange3 = $d.Range($byPara.Range.Start, $byPara.Range.End)
$byRange3.Find.Execute("),") | Out-Null
$byRange3.Collapse(0)               # collapse to position right after "),"
$byRange3.InsertParagraphAfter()

# The new paragraph now starts with " pxr180025." -- remove the leading
# space and prepend "Prajakta Ray (".
$prajPara = $d.Paragraphs.Item(3)
$leadSpace = $d.Range($prajPara.Range.Start, $prajPara.Range.Start + 1)
$leadSpace.Delete()

$prajPara2 = $d.Paragraphs.Item(3)
$prajStart = $d.Range($prajPara2.Range.Start, $prajPara2.Range.Start)
$prajStart.InsertBefore("Prajakta Ray (")

# Insert ")" right before the final "." of this paragraph.
$prajPara3 = $d.Paragraphs.Item(3)
$prajRange = $d.Range($prajPara3.Range.Start, $prajPara3.Range.End)
$prajRange.Find.Execute("pxr180025") | Out-Null
$prajRange.Collapse(0)
$prajRange.InsertBefore(")")

# ------------------------------------------------------------------
# Step 3: update the "Readme" paragraph to "Readme – Bounded Queue"
# and re-create the "_GoBack" bookmark at the very end of it.
# ------------------------------------------------------------------
$readmePara = $d.Paragraphs.Item(5)
$readmeRange = $d.Range($readmePara.Range.Start, $readmePara.Range.End)
$readmeRange.Find.Execute("Readme") | Out-Null
$readmeRange.Collapse(0)
$readmeRange.InsertBefore(" – Bounded Queue")

$readmePara2 = $d.Paragraphs.Item(5)
$readmeEnd = $readmePara2.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($readmeEnd, $readmeEnd))
